# Generate Report for Handoff
# Adds a new file's rows (8a3b0a26-71be-4583-a23f-19f68e58a3fe) to the
# "Overview", "zh-cn" and "de-de" tabs, following the same layout as the
# existing 26bc3644-... entry.

$wb = $excel.ActiveWorkbook

$guid = "8a3b0a26-71be-4583-a23f-19f68e58a3fe"
$zhHash = "ef41eb1ecace52797d77505df22fcb13f411ec9b"

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$tblOverview = $wsOverview.ListObjects.Item("Overview")
$rowOverview = $tblOverview.ListRows.Add()
$rngOverview = $rowOverview.Range

$rngOverview.Cells.Item(1, 1).Value = "$guid.md"
$rngOverview.Cells.Item(1, 2).Value = "e2e\$guid.md"
$rngOverview.Cells.Item(1, 3).Value = ".md"
$rngOverview.Cells.Item(1, 4).Value = ""
$rngOverview.Cells.Item(1, 5).Value = "Ready for handoff"
$rngOverview.Cells.Item(1, 6).Value = "Ready for handoff"
$rngOverview.Cells.Item(1, 7).Value = "2016-08-21 20:51:47"

$wsOverview.Hyperlinks.Add(
    $rngOverview.Cells.Item(1, 2),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8ec94dd775a68cd5cc747e1065ee4baf22388889/e2e/$guid.md",
    "",
    "",
    "e2e\$guid.md"
)

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$tblZhCn = $wsZhCn.ListObjects.Item("zh-cn")
$rowZhCn = $tblZhCn.ListRows.Add()
$rngZhCn = $rowZhCn.Range

$rngZhCn.Cells.Item(1, 1).Value = "$guid.md"
$rngZhCn.Cells.Item(1, 2).Value = ".md"
$rngZhCn.Cells.Item(1, 3).Value = "Ready for handoff"
$rngZhCn.Cells.Item(1, 4).Value = "e2e"
$rngZhCn.Cells.Item(1, 5).Value = "ht"
$rngZhCn.Cells.Item(1, 6).Value = "False"
$rngZhCn.Cells.Item(1, 7).Value = "$guid.$zhHash.zh-cn.xlf"
$rngZhCn.Cells.Item(1, 8).Value = "2016-08-21 20:51:43"
$rngZhCn.Cells.Item(1, 9).Value = ""
$rngZhCn.Cells.Item(1, 10).Value = ""
$rngZhCn.Cells.Item(1, 11).Value = "0001-01-01 00:00:00"
$rngZhCn.Cells.Item(1, 12).Value = ""
$rngZhCn.Cells.Item(1, 13).Value = "True"
$rngZhCn.Cells.Item(1, 14).Value = ""
$rngZhCn.Cells.Item(1, 15).Value = "False"
$rngZhCn.Cells.Item(1, 16).Value = ""

$wsZhCn.Hyperlinks.Add(
    $rngZhCn.Cells.Item(1, 1),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8ec94dd775a68cd5cc747e1065ee4baf22388889/e2e/$guid.md",
    "",
    "",
    "$guid.md"
)

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$tblDeDe = $wsDeDe.ListObjects.Item("de-de")
$rowDeDe = $tblDeDe.ListRows.Add()
$rngDeDe = $rowDeDe.Range

$rngDeDe.Cells.Item(1, 1).Value = "$guid.md"
$rngDeDe.Cells.Item(1, 2).Value = ".md"
$rngDeDe.Cells.Item(1, 3).Value = "Ready for handoff"
$rngDeDe.Cells.Item(1, 4).Value = "e2e"
$rngDeDe.Cells.Item(1, 5).Value = "ht"
$rngDeDe.Cells.Item(1, 6).Value = "False"
$rngDeDe.Cells.Item(1, 7).Value = "$guid.$zhHash.de-de.xlf"
$rngDeDe.Cells.Item(1, 8).Value = "2016-08-21 20:51:47"
$rngDeDe.Cells.Item(1, 9).Value = ""
$rngDeDe.Cells.Item(1, 10).Value = ""
$rngDeDe.Cells.Item(1, 11).Value = "0001-01-01 00:00:00"
$rngDeDe.Cells.Item(1, 12).Value = ""
$rngDeDe.Cells.Item(1, 13).Value = "True"
$rngDeDe.Cells.Item(1, 14).Value = ""
$rngDeDe.Cells.Item(1, 15).Value = "False"
$rngDeDe.Cells.Item(1, 16).Value = ""

$wsDeDe.Hyperlinks.Add(
    $rngDeDe.Cells.Item(1, 1),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8ec94dd775a68cd5cc747e1065ee4baf22388889/e2e/$guid.md",
    "",
    "",
    "$guid.md"
)
